# Added two new Mac-Addresses (new rows of reg_center_user_machine test data)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @{ Row = 31; A = 10001; B = 110030; C = 10030 },
    @{ Row = 32; A = 10001; B = 110031; C = 10031 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = "eng"
    $ws.Range("E$row").Value = $true
    $ws.Range("F$row").Value = "superadmin"
    $ws.Range("G$row").Value = "now()"
    $ws.Range("H$row").Value = "now()"
}

# Mirror the author's final selection/scroll position after entering the new rows
$ws.Range("F30").Select()
